# Update the "想去人数" (interested-people count) figures in column F
# on the "展览" sheet and the "全部类型" sheet, matching the refreshed
# output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 256
$ws1.Range("F3").Value  = 1093
$ws1.Range("F5").Value  = 441
$ws1.Range("F6").Value  = 83
$ws1.Range("F7").Value  = 572
$ws1.Range("F9").Value  = 6867
$ws1.Range("F15").Value = 1115
$ws1.Range("F16").Value = 16312
$ws1.Range("F19").Value = 337
$ws1.Range("F20").Value = 190
$ws1.Range("F22").Value = 11434
$ws1.Range("F23").Value = 14
$ws1.Range("F24").Value = 1071
$ws1.Range("F25").Value = 4502
$ws1.Range("F26").Value = 358
$ws1.Range("F27").Value = 390
$ws1.Range("F30").Value = 321
$ws1.Range("F32").Value = 5214

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 256
$ws4.Range("F3").Value  = 1093
$ws4.Range("F5").Value  = 441
$ws4.Range("F6").Value  = 83
$ws4.Range("F7").Value  = 572
$ws4.Range("F10").Value = 6867
$ws4.Range("F17").Value = 1115
$ws4.Range("F18").Value = 16312
$ws4.Range("F21").Value = 337
$ws4.Range("F22").Value = 190
$ws4.Range("F26").Value = 11434
$ws4.Range("F27").Value = 14
$ws4.Range("F28").Value = 1071
$ws4.Range("F29").Value = 4502
$ws4.Range("F30").Value = 358
$ws4.Range("F31").Value = 390
$ws4.Range("F34").Value = 321
$ws4.Range("F36").Value = 5214
